# Lab Exam 03 grading pass: fill in the "Total Points" (column E) scores
# for the Customer Class and Product Class sections, mirroring the
# "Points for grading" (column D) values recorded by the grader.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer Class section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Product Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the view where the grader last clicked, at the running total cell
$ws.Range("E15").Select() | Out-Null
